$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction: the I and J columns held values that were a factor of 10 too
# small (decimal point slipped during the original entry) - multiply them
# back up for rows 2-22.
for ($r = 2; $r -le 22; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2
    $jVal = $ws.Cells.Item($r, 10).Value2
    $ws.Cells.Item($r, 9).Value = $iVal * 10
    $ws.Cells.Item($r, 10).Value = $jVal * 10
}

# Columns L and M get a fixed width of 10
$ws.Columns("L:M").ColumnWidth = 9.166666666666666

# Move the active selection
$ws.Range("O15").Select() | Out-Null

Write-Host "done"
